$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values that look numeric are written as text (to preserve
# the "x.xxx.xx"-style display strings exactly as in the source data).
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "42.356.66"
$ws.Cells.Item(2, 5).Value = "  -2.86%  "
$ws.Cells.Item(3, 4).Value = "2.223.51"
$ws.Cells.Item(3, 5).Value = "  -1.97%  "
$ws.Cells.Item(4, 5).Value = "  +0.02%  "
$ws.Cells.Item(5, 4).Value = "110.69"
$ws.Cells.Item(5, 5).Value = "  -8.29%  "
$ws.Cells.Item(6, 4).Value = "290.35"
$ws.Cells.Item(6, 5).Value = "  +7.79%  "
$ws.Cells.Item(7, 4).Value = "0.623"
$ws.Cells.Item(7, 5).Value = "  -2.80%  "
$ws.Cells.Item(8, 5).Value = "  -0.39%  "
$ws.Cells.Item(9, 4).Value = "0.599"
$ws.Cells.Item(9, 5).Value = "  -3.61%  "
$ws.Cells.Item(10, 4).Value = "43.51"
$ws.Cells.Item(10, 5).Value = "  -8.85%  "
$ws.Cells.Item(11, 4).Value = "0.0909"
$ws.Cells.Item(11, 5).Value = "  -3.72%  "
$ws.Cells.Item(12, 4).Value = "54.12"
$ws.Cells.Item(13, 4).Value = "8.60"
$ws.Cells.Item(13, 5).Value = "  -8.92%  "
$ws.Cells.Item(14, 4).Value = "1.01"
$ws.Cells.Item(14, 5).Value = "  +10.09%  "
$ws.Cells.Item(15, 5).Value = "  -3.01%  "
$ws.Cells.Item(16, 4).Value = "14.87"
$ws.Cells.Item(16, 5).Value = "  -6.32%  "
$ws.Cells.Item(17, 4).Value = "2.559.34"
$ws.Cells.Item(17, 5).Value = "  -1.91%  "
$ws.Cells.Item(18, 4).Value = "2.220.54"
$ws.Cells.Item(18, 5).Value = "  -1.95%  "
$ws.Cells.Item(19, 4).Value = "42.349.03"
$ws.Cells.Item(19, 5).Value = "  -2.85%  "
$ws.Cells.Item(20, 4).Value = "7.11"
$ws.Cells.Item(20, 5).Value = "  +2.83%  "
$ws.Cells.Item(21, 5).Value = "  -4.70%  "
$ws.Cells.Item(22, 4).Value = "72.68"
$ws.Cells.Item(22, 5).Value = "  -0.12%  "
$ws.Cells.Item(23, 5).Value = "  +13.12%  "
$ws.Cells.Item(24, 4).Value = "2.38"
$ws.Cells.Item(24, 5).Value = "  -0.48%  "
$ws.Cells.Item(25, 4).Value = "229.76"
$ws.Cells.Item(25, 5).Value = "  -2.23%  "
$ws.Cells.Item(26, 4).Value = "8.97"
$ws.Cells.Item(26, 5).Value = "  -6.76%  "
$ws.Cells.Item(27, 4).Value = "1.00"
$ws.Cells.Item(27, 5).Value = "  -1.77%  "
$ws.Cells.Item(28, 4).Value = "11.40"
$ws.Cells.Item(28, 5).Value = "  -6.86%  "
$ws.Cells.Item(29, 5).Value = "  -2.64%  "
$ws.Cells.Item(30, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(30, 4).Value = "37.43"
$ws.Cells.Item(30, 5).Value = "  -11.45%  "
$ws.Cells.Item(31, 2).Value = "Monero"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(31, 4).Value = "173.35"
$ws.Cells.Item(31, 5).Value = "  -0.93%  "
$ws.Cells.Item(32, 4).Value = "3.08"
$ws.Cells.Item(32, 5).Value = "  -8.27%  "
$ws.Cells.Item(33, 4).Value = "20.82"
$ws.Cells.Item(33, 5).Value = "  -3.35%  "
$ws.Cells.Item(34, 4).Value = "0.0880"
$ws.Cells.Item(34, 5).Value = "  -4.22%  "
$ws.Cells.Item(35, 4).Value = "5.60"
$ws.Cells.Item(35, 5).Value = "  -2.01%  "
$ws.Cells.Item(36, 4).Value = "4.97"
$ws.Cells.Item(36, 5).Value = "  +5.32%  "
$ws.Cells.Item(37, 4).Value = "0.126"
$ws.Cells.Item(37, 5).Value = "  -3.55%  "
$ws.Cells.Item(38, 4).Value = "4.17"
$ws.Cells.Item(38, 5).Value = "  -7.07%  "
$ws.Cells.Item(39, 4).Value = "0.0371"
$ws.Cells.Item(39, 5).Value = "  -2.37%  "
$ws.Cells.Item(40, 4).Value = "0.104"
$ws.Cells.Item(40, 5).Value = "  -4.87%  "
$ws.Cells.Item(41, 2).Value = "MultiversX"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Cells.Item(41, 4).Value = "73.70"
$ws.Cells.Item(41, 5).Value = "  +1.82%  "
$ws.Cells.Item(42, 2).Value = "LidoDAOToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(42, 4).Value = "2.39"
$ws.Cells.Item(42, 5).Value = "  -6.37%  "
$ws.Cells.Item(43, 4).Value = "0.231"
$ws.Cells.Item(43, 5).Value = "  -4.27%  "
$ws.Cells.Item(44, 4).Value = "0.999"
$ws.Cells.Item(44, 5).Value = "  -0.23%  "
$ws.Cells.Item(45, 4).Value = "12.26"
$ws.Cells.Item(45, 5).Value = "  -10.91%  "
$ws.Cells.Item(46, 4).Value = "1.30"
$ws.Cells.Item(46, 5).Value = "  -4.91%  "
$ws.Cells.Item(47, 4).Value = "5.32"
$ws.Cells.Item(47, 5).Value = "  -6.76%  "
$ws.Cells.Item(48, 4).Value = "1.71"
$ws.Cells.Item(48, 5).Value = "  +6.99%  "
$ws.Cells.Item(49, 4).Value = "1.27"
$ws.Cells.Item(49, 5).Value = "  -0.67%  "
$ws.Cells.Item(50, 2).Value = "Aave"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(50, 4).Value = "100.94"
$ws.Cells.Item(50, 5).Value = "  -1.86%  "
$ws.Cells.Item(51, 2).Value = "FraxShare"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(51, 4).Value = "8.39"
$ws.Cells.Item(51, 5).Value = "  -2.20%  "

# Reset style on column D so no stray style index is introduced (keeps cells
# attribute-for-attribute consistent with the rest of the unstyled data cells).
$dRange.Style = "Normal"

